$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sponza" (2nd sheet) - add new column L (v1423) after existing K
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Copy the formatting (styles) of column K into column L first, so every
# cell in the new column carries the same style index as its K neighbour.
$ws2.Range("K1:K16").Copy() | Out-Null
$ws2.Range("L1:L16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Header
$ws2.Range("L1").Value = "v1423"

# Raw sample data
$ws2.Range("L2").Value = 7242
$ws2.Range("L3").Value = 7229
$ws2.Range("L4").Value = 7226
$ws2.Range("L5").Value = 7254
$ws2.Range("L6").Value = 7204
$ws2.Range("L7").Value = 7226
$ws2.Range("L8").Value = 7212
$ws2.Range("L9").Value = 7262
$ws2.Range("L10").Value = 7229
$ws2.Range("L11").Value = 7264

# Summary rows
$ws2.Range("L12").Formula = "=AVERAGE(L2:L11)"
$ws2.Range("L13").Formula = "=_xlfn.VAR.S(L2:L11)"
$ws2.Range("L14").Formula = "=1-_xlfn.T.TEST(K2:K11,L2:L11,2,3)"
$ws2.Range("L15").Formula = "=B12/L12"
$ws2.Range("L16").Formula = "=B12/L12"

# Extend the conditional formatting that highlights the ratio rows
$fcs2 = $ws2.Range("B15:K16").FormatConditions
for ($i = 1; $i -le $fcs2.Count; $i++) {
    $fcs2.Item($i).ModifyAppliesToRange($ws2.Range("B15:L16"))
}

# Update selection to match the authored workbook
$ws2.Activate()
$ws2.Range("L2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "ComplexMesh" (3rd sheet) - add new column K (v1423) after existing J
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Copy the formatting (styles) of column J into column K first.
$ws3.Range("J1:J16").Copy() | Out-Null
$ws3.Range("K1:K16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Header
$ws3.Range("K1").Value = "v1423"

# Raw sample data
$ws3.Range("K2").Value = 5249
$ws3.Range("K3").Value = 5080
$ws3.Range("K4").Value = 5073
$ws3.Range("K5").Value = 5055
$ws3.Range("K6").Value = 5051
$ws3.Range("K7").Value = 5085
$ws3.Range("K8").Value = 5050
$ws3.Range("K9").Value = 5044
$ws3.Range("K10").Value = 5058
$ws3.Range("K11").Value = 5077

# Summary rows (K12 is a literal value here, mirroring H12/I12/J12 on this sheet)
$ws3.Range("K12").Value = 5119
$ws3.Range("K13").Formula = "=_xlfn.VAR.S(K2:K11)"
$ws3.Range("K14").Formula = "=1-_xlfn.T.TEST(J2:J11,K2:K11,2,3)"
$ws3.Range("K15").Formula = "=B12/K12"
$ws3.Range("K16").Formula = "=B12/K12"

# Extend the conditional formatting that highlights the ratio rows
$fcs3 = $ws3.Range("B15:J16").FormatConditions
for ($i = 1; $i -le $fcs3.Count; $i++) {
    $fcs3.Item($i).ModifyAppliesToRange($ws3.Range("B15:K16"))
}

# Update selection to match the authored workbook
$ws3.Activate()
$ws3.Range("K4").Select() | Out-Null
